$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Precinct crime statistics table updates (rows 15-30, 33) ---
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 100
$ws.Cells.Item(15, 6).Value = 5
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 8).Value = 66.666666666666
$ws.Cells.Item(15, 9).Value = 26
$ws.Cells.Item(15, 10).Value = 19
$ws.Cells.Item(15, 11).Value = 36.842105263157
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 188.888888888889
$ws.Cells.Item(15, 14).Value = -16.129032258064
$ws.Cells.Item(16, 3).Value = 16
$ws.Cells.Item(16, 4).Value = 11
$ws.Cells.Item(16, 5).Value = 45.454545454545
$ws.Cells.Item(16, 6).Value = 63
$ws.Cells.Item(16, 7).Value = 52
$ws.Cells.Item(16, 8).Value = 21.153846153846
$ws.Cells.Item(16, 9).Value = 269
$ws.Cells.Item(16, 10).Value = 277
$ws.Cells.Item(16, 11).Value = -2.888086642599
$ws.Cells.Item(16, 12).Value = 4.263565891472
$ws.Cells.Item(16, 13).Value = 31.862745098039
$ws.Cells.Item(16, 14).Value = -67.97619047619
$ws.Cells.Item(17, 4).Value = 21
$ws.Cells.Item(17, 5).Value = -4.761904761904
$ws.Cells.Item(17, 7).Value = 83
$ws.Cells.Item(17, 8).Value = -2.409638554216
$ws.Cells.Item(17, 9).Value = 455
$ws.Cells.Item(17, 10).Value = 449
$ws.Cells.Item(17, 11).Value = 1.336302895322
$ws.Cells.Item(17, 12).Value = 6.55737704918
$ws.Cells.Item(17, 13).Value = 121.951219512195
$ws.Cells.Item(17, 14).Value = -2.985074626865
$ws.Cells.Item(18, 3).Value = 10
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 5).Value = 150
$ws.Cells.Item(18, 7).Value = 32
$ws.Cells.Item(18, 8).Value = 12.5
$ws.Cells.Item(18, 9).Value = 188
$ws.Cells.Item(18, 10).Value = 164
$ws.Cells.Item(18, 11).Value = 14.634146341463
$ws.Cells.Item(18, 12).Value = 35.251798561151
$ws.Cells.Item(18, 13).Value = 97.894736842105
$ws.Cells.Item(18, 14).Value = -69.96805111821
$ws.Cells.Item(19, 3).Value = 16
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 77.777777777777
$ws.Cells.Item(19, 6).Value = 73
$ws.Cells.Item(19, 7).Value = 57
$ws.Cells.Item(19, 8).Value = 28.070175438596
$ws.Cells.Item(19, 9).Value = 413
$ws.Cells.Item(19, 10).Value = 404
$ws.Cells.Item(19, 11).Value = 2.227722772277
$ws.Cells.Item(19, 12).Value = 33.656957928802
$ws.Cells.Item(19, 13).Value = 141.520467836257
$ws.Cells.Item(19, 14).Value = 31.528662420382
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 4).Value = 9
$ws.Cells.Item(20, 5).Value = -33.333333333333
$ws.Cells.Item(20, 6).Value = 19
$ws.Cells.Item(20, 7).Value = 18
$ws.Cells.Item(20, 8).Value = 5.555555555555
$ws.Cells.Item(20, 9).Value = 119
$ws.Cells.Item(20, 10).Value = 94
$ws.Cells.Item(20, 11).Value = 26.595744680851
$ws.Cells.Item(20, 12).Value = -24.683544303797
$ws.Cells.Item(20, 13).Value = 128.846153846154
$ws.Cells.Item(20, 14).Value = -59.523809523809
$ws.Cells.Item(21, 3).Value = 70
$ws.Cells.Item(21, 4).Value = 55
$ws.Cells.Item(21, 5).Value = 27.272727272727
$ws.Cells.Item(21, 6).Value = 277
$ws.Cells.Item(21, 8).Value = 13.061224489795
$ws.Cells.Item(21, 9).Value = 1476
$ws.Cells.Item(21, 10).Value = 1413
$ws.Cells.Item(21, 11).Value = 4.458598726114
$ws.Cells.Item(21, 12).Value = 11.649016641452
$ws.Cells.Item(21, 13).Value = 99.459459459459
$ws.Cells.Item(21, 14).Value = -43.40490797546
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 3).NumberFormat = "#,##0"
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 7).Value = 8
$ws.Cells.Item(22, 8).Value = -50
$ws.Cells.Item(22, 9).Value = 31
$ws.Cells.Item(22, 10).Value = 37
$ws.Cells.Item(22, 11).Value = -16.216216216216
$ws.Cells.Item(22, 12).Value = 14.814814814814
$ws.Cells.Item(22, 13).Value = 10.714285714285
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 4).Value = 9
$ws.Cells.Item(23, 5).Value = -33.333333333333
$ws.Cells.Item(23, 6).Value = 33
$ws.Cells.Item(23, 8).Value = 22.222222222222
$ws.Cells.Item(23, 9).Value = 223
$ws.Cells.Item(23, 10).Value = 216
$ws.Cells.Item(23, 11).Value = 3.24074074074
$ws.Cells.Item(23, 12).Value = -7.851239669421
$ws.Cells.Item(23, 13).Value = 72.868217054263
$ws.Cells.Item(24, 3).Value = 26
$ws.Cells.Item(24, 4).Value = 21
$ws.Cells.Item(24, 5).Value = 23.809523809523
$ws.Cells.Item(24, 6).Value = 127
$ws.Cells.Item(24, 7).Value = 81
$ws.Cells.Item(24, 8).Value = 56.79012345679
$ws.Cells.Item(24, 9).Value = 782
$ws.Cells.Item(24, 10).Value = 689
$ws.Cells.Item(24, 11).Value = 13.497822931785
$ws.Cells.Item(24, 12).Value = 11.237553342816
$ws.Cells.Item(24, 13).Value = 36.713286713286
$ws.Cells.Item(25, 3).Value = 16
$ws.Cells.Item(25, 4).Value = 9
$ws.Cells.Item(25, 5).Value = 77.777777777777
$ws.Cells.Item(25, 6).Value = 61
$ws.Cells.Item(25, 7).Value = 28
$ws.Cells.Item(25, 8).Value = 117.857142857143
$ws.Cells.Item(25, 9).Value = 317
$ws.Cells.Item(25, 10).Value = 277
$ws.Cells.Item(25, 11).Value = 14.440433212996
$ws.Cells.Item(25, 12).Value = 20.532319391635
$ws.Cells.Item(26, 3).Value = 20
$ws.Cells.Item(26, 4).Value = 33
$ws.Cells.Item(26, 5).Value = -39.393939393939
$ws.Cells.Item(26, 6).Value = 90
$ws.Cells.Item(26, 7).Value = 113
$ws.Cells.Item(26, 8).Value = -20.353982300885
$ws.Cells.Item(26, 9).Value = 563
$ws.Cells.Item(26, 10).Value = 564
$ws.Cells.Item(26, 11).Value = -0.177304964539
$ws.Cells.Item(26, 12).Value = 12.15139442231
$ws.Cells.Item(26, 13).Value = 13.052208835341
$ws.Cells.Item(27, 3).Value = 3
$ws.Cells.Item(27, 5).Value = 200
$ws.Cells.Item(27, 6).Value = 8
$ws.Cells.Item(27, 7).Value = 3
$ws.Cells.Item(27, 8).Value = 166.666666666667
$ws.Cells.Item(27, 9).Value = 34
$ws.Cells.Item(27, 10).Value = 26
$ws.Cells.Item(27, 11).Value = 30.76923076923
$ws.Cells.Item(27, 12).Value = -5.555555555555
$ws.Cells.Item(28, 3).Value = 4
$ws.Cells.Item(28, 4).Value = 4
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 10
$ws.Cells.Item(28, 7).Value = 14
$ws.Cells.Item(28, 8).Value = -28.571428571428
$ws.Cells.Item(28, 9).Value = 53
$ws.Cells.Item(28, 10).Value = 75
$ws.Cells.Item(28, 11).Value = -29.333333333333
$ws.Cells.Item(28, 12).Value = -10.169491525423
$ws.Cells.Item(29, 10).Value = 16
$ws.Cells.Item(29, 11).Value = 25
$ws.Cells.Item(29, 12).Value = 81.818181818181
$ws.Cells.Item(29, 13).Value = -20
$ws.Cells.Item(29, 14).Value = -77.011494252873
$ws.Cells.Item(30, 10).Value = 15
$ws.Cells.Item(30, 11).Value = 13.333333333333
$ws.Cells.Item(30, 12).Value = 54.545454545454
$ws.Cells.Item(30, 13).Value = -5.555555555555
$ws.Cells.Item(30, 14).Value = -77.922077922077
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 3).NumberFormat = "#,##0"
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 6).NumberFormat = "#,##0"
$ws.Cells.Item(33, 9).Value = 3
$ws.Cells.Item(33, 11).Value = 200
$ws.Cells.Item(33, 12).Value = 0
